# "add preview in datasets"
# Inserts two new metadata rows (dataset.preview.table / dataset.preview.line)
# right after the existing "dataset.note" row (i.e. before the old row 4),
# pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert two blank rows at row 4 - everything from the old row 4 onward
# (dataset.note, dataset.source, ...) shifts down to rows 6+.
$ws.Rows.Item(4).Resize(2).Insert()

# New row 4: dataset.preview.table key/value pair.
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"

# New row 5: dataset.preview.line key/value pair.
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

# Match formatting of the new cells: wrapped, vertically centred, tall rows
# so the multi-line formula text is fully visible.
$ws.Range("A4:B5").VerticalAlignment = -4108
$ws.Range("A4:B5").WrapText = $true
$ws.Range("A4:B5").RowHeight = 120

# Leave the selection where the author left it after typing the new rows.
$ws.Range("B8").Select()
